# Timesheet Calculator v2 - fix credited/excess minute computation for the
# last working day of the 3rd week (19-Jun-2021, row 22) and append the
# missing 4th "WEEK COVERED" block (21-Jun-2021) to the employee sheet,
# mirroring the existing week blocks. Also restores the natural
# sheet/selection bookkeeping (the employee sheet becomes the active tab).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$emp     = $wb.Worksheets.Item("Doe, Jean S.")

# --- Bug fix: 19-Jun-2021 (row 22) had its credited/excess minutes swapped ---
$emp.Range("E22").Value = 0
$emp.Range("F22").Value = 120

# --- New week block: WEEK COVERED: 21-Jun-2021 to 21-Jun-2021 (row 25) ---
$weekHeader = $emp.Range("A25")
$weekHeader.Value = "WEEK COVERED: 21-Jun-2021 to 21-Jun-2021"
$weekHeader.Font.Name = "Calibri"
$weekHeader.Font.Size = 11

# --- Column headers for the new week (row 26) ---
$emp.Range("A26").Value = "Date"
$emp.Range("B26").Value = "Day"
$emp.Range("C26").Value = "Time-in/Time-out"
$emp.Range("D26").Value = "Rendered MINS for the Day"
$emp.Range("E26").Value = "Credited Regular Log [480 = 1 day]"
$emp.Range("F26").Value = "Minutes in excess of 480; Sat/Sun Duties"

# --- Single day row for 21-Jun-2021, a Monday with no logged time (row 27) ---
$dateCell = $emp.Range("A27")
$dateCell.NumberFormat = "@"
$dateCell.Value = "21-Jun-2021"
$dateCell.NumberFormat = "General"

$emp.Range("B27").Value = "MON"
$emp.Range("D27").Value = 0
$emp.Range("E27").Value = 0
$emp.Range("F27").Value = 0

# --- Selection / active sheet bookkeeping ---
$summary.Range("A3").Select() | Out-Null
$emp.Activate() | Out-Null
$emp.Range("D33").Select() | Out-Null
